# Update the "Förändrad" (changed) date column C for rows 2-31
# from 45207 (2023-10-08) to 45208 (2023-10-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}
